$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 5882911.5
$ws.Range("I92").Value = 540.5
$ws.Range("K92").Value = 540.5
$ws.Range("M92").Value = 707.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 2330.2222
$ws.Range("I129").Value = 1287
$ws.Range("J129").Value = 2943.8823
$ws.Range("K129").Value = 3861
$ws.Range("L129").Value = 8831.6469
$ws.Range("M129").Value = 1139
$ws.Range("N129").Value = -18831.6469

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H131").Value = 589581.5
$ws.Range("I131").Value = 589581.5
$ws.Range("K131").Value = 1768744.5
$ws.Range("M131").Value = -1763704.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 4695.706
$ws.Range("I132").Value = 5123.6
$ws.Range("K132").Value = 15370.8
$ws.Range("M132").Value = -12840.8

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2285.5293
$ws.Range("I137").Value = 1859.8422
$ws.Range("K137").Value = 5579.5266
$ws.Range("M137").Value = -3029.5266

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2943.52
$ws.Range("I138").Value = 2415.4614
$ws.Range("J138").Value = 3129.054
$ws.Range("K138").Value = 7246.3842
$ws.Range("L138").Value = 9387.162
$ws.Range("M138").Value = -2106.3842
$ws.Range("N138").Value = -19667.162

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2235.63
$ws.Range("I32").Value = 1840.8866
$ws.Range("K32").Value = 1840.8866
$ws.Range("M32").Value = -1553.8866

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 8654.833000000001
$ws.Range("I61").Value = 8027.4287
$ws.Range("K61").Value = 8027.4287
$ws.Range("M61").Value = -7815.4287

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 141465.8
$ws.Range("I74").Value = 187320.73
$ws.Range("K74").Value = 187320.73
$ws.Range("M74").Value = -186446.73

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 141465.8
$ws.Range("I77").Value = 187320.73
$ws.Range("K77").Value = 936603.65
$ws.Range("M77").Value = -932235.65

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 8654.833000000001
$ws.Range("I136").Value = 8027.4287
$ws.Range("K136").Value = 24082.2861
$ws.Range("M136").Value = -21532.2861

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 30303566
$ws.Range("I94").Value = 35714730
$ws.Range("J94").Value = 1069.2
$ws.Range("K94").Value = 35714730
$ws.Range("L94").Value = 1069.2
$ws.Range("M94").Value = -35714279
$ws.Range("N94").Value = -1971.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1717.6364
$ws.Range("I107").Value = 1539.4
$ws.Range("J107").Value = 3500
$ws.Range("K107").Value = 1539.4
$ws.Range("L107").Value = 3500
$ws.Range("M107").Value = 380.5999999999999
$ws.Range("N107").Value = -7340

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2357.7144
$ws.Range("I134").Value = 2100.68
$ws.Range("J134").Value = 4499.6665
$ws.Range("K134").Value = 6302.039999999999
$ws.Range("L134").Value = 13498.9995
$ws.Range("M134").Value = -3767.039999999999
$ws.Range("N134").Value = -18568.9995

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3041.7273
$ws.Range("I16").Value = 2529.25
$ws.Range("K16").Value = 2529.25
$ws.Range("M16").Value = -2242.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3124.774
$ws.Range("I31").Value = 2635.9539
$ws.Range("J31").Value = 4797.0527
$ws.Range("K31").Value = 2635.9539
$ws.Range("L31").Value = 4797.0527
$ws.Range("M31").Value = -2340.9539
$ws.Range("N31").Value = -5387.0527

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 3124.774
$ws.Range("I34").Value = 2635.9539
$ws.Range("J34").Value = 4797.0527
$ws.Range("K34").Value = 2635.9539
$ws.Range("L34").Value = 4797.0527
$ws.Range("M34").Value = -2433.9539
$ws.Range("N34").Value = -5201.0527

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3346.682
$ws.Range("I58").Value = 2354.5454
$ws.Range("J58").Value = 4338.8184
$ws.Range("K58").Value = 2354.5454
$ws.Range("L58").Value = 4338.8184
$ws.Range("M58").Value = -2151.5454
$ws.Range("N58").Value = -4744.8184

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H88").Value = 41996
$ws.Range("J88").Value = 41996
$ws.Range("L88").Value = 41996
$ws.Range("N88").Value = -42808

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H91").Value = 41996
$ws.Range("J91").Value = 41996
$ws.Range("L91").Value = 41996
$ws.Range("N91").Value = -44804

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 736.5217
$ws.Range("I107").Value = 752.3889
$ws.Range("J107").Value = 679.4
$ws.Range("K107").Value = 752.3889
$ws.Range("L107").Value = 679.4
$ws.Range("M107").Value = 1167.6111
$ws.Range("N107").Value = -4519.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 3041.7273
$ws.Range("I113").Value = 2529.25
$ws.Range("K113").Value = 2529.25
$ws.Range("M113").Value = -359.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 3346.682
$ws.Range("I136").Value = 2354.5454
$ws.Range("J136").Value = 4338.8184
$ws.Range("K136").Value = 7063.6362
$ws.Range("L136").Value = 13016.4552
$ws.Range("M136").Value = -4513.6362
$ws.Range("N136").Value = -18116.4552

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 8506.75
$ws.Range("I3").Value = 8649.166999999999
$ws.Range("K3").Value = 25947.501
$ws.Range("M3").Value = -25835.501

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 11243.5
$ws.Range("J39").Value = 11243.5
$ws.Range("L39").Value = 33730.5
$ws.Range("N39").Value = -34318.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 1040298.4
$ws.Range("J97").Value = 100248.5
$ws.Range("L97").Value = 300745.5
$ws.Range("N97").Value = -301737.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 400919.6
$ws.Range("I121").Value = 530
$ws.Range("J121").Value = 501017
$ws.Range("K121").Value = 1590
$ws.Range("L121").Value = 1503051
$ws.Range("M121").Value = -280
$ws.Range("N121").Value = -1505671

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 2294.8572
$ws.Range("J129").Value = 2498.4
$ws.Range("L129").Value = 7495.200000000001
$ws.Range("N129").Value = -17495.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 6192.6763
$ws.Range("I131").Value = 25505
$ws.Range("J131").Value = 2054.3215
$ws.Range("K131").Value = 76515
$ws.Range("L131").Value = 6162.9645
$ws.Range("M131").Value = -71475
$ws.Range("N131").Value = -16242.9645

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H139").Value = 3626.8125
$ws.Range("I139").Value = 3147
$ws.Range("J139").Value = 4000
$ws.Range("K139").Value = 9441
$ws.Range("L139").Value = 12000
$ws.Range("M139").Value = -4301
$ws.Range("N139").Value = -22280

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2471.9167
$ws.Range("J122").Value = 1799
$ws.Range("L122").Value = 5397
$ws.Range("N122").Value = -10297

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 5981.364
$ws.Range("I132").Value = 4828.7144
$ws.Range("J132").Value = 7998.5
$ws.Range("K132").Value = 14486.1432
$ws.Range("L132").Value = 23995.5
$ws.Range("M132").Value = -11956.1432
$ws.Range("N132").Value = -29055.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1298
$ws.Range("I16").Value = 1298
$ws.Range("K16").Value = 1298
$ws.Range("M16").Value = -1128

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H23").Value = 4999
$ws.Range("I23").Value = 4999
$ws.Range("K23").Value = 4999
$ws.Range("M23").Value = -4769

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 21289.07
$ws.Range("I40").Value = 35633.094
$ws.Range("K40").Value = 35633.094
$ws.Range("M40").Value = -35497.094

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4524.9
$ws.Range("I132").Value = 4035.5715
$ws.Range("J132").Value = 5666.6665
$ws.Range("K132").Value = 12106.7145
$ws.Range("L132").Value = 16999.9995
$ws.Range("M132").Value = -9576.7145
$ws.Range("N132").Value = -22059.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H50").Value = 9874.25
$ws.Range("J50").Value = 9874.25
$ws.Range("L50").Value = 9874.25
$ws.Range("N50").Value = -11136.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3936.9285
$ws.Range("I132").Value = 3621.8164
$ws.Range("K132").Value = 10865.4492
$ws.Range("M132").Value = -8335.449200000001
